# feat: add 2022-Q4 data
#
# Strategy (mirrors how the real diff landed, so sheetIds/rIds line up):
#   1. Duplicate the existing "2022-Q3" sheet (sheet index 2) right after
#      itself. The duplicate keeps an exact copy of the old per-fund
#      breakdown, formatting and all.
#   2. Rename the *original* sheet to "2022-Q4" and overwrite its cells
#      in place with the new quarter's fund table - this sheet keeps its
#      original sheetId/relationship id, just like the source edit.
#   3. Rename the *duplicate* to "2022-Q3" - it already holds the old
#      quarter's untouched data, so nothing else to do there.
#   4. Update the "总计" (totals) summary sheet: insert the 2022-Q4 totals
#      as the new row 2 and push the old 2022-Q3 totals down to row 3.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force Excel to store the value as literal text (inlineStr) even when
    # it looks numeric (e.g. "012093" or "5.64"), the same way a leading
    # apostrophe does in the Excel UI. The apostrophe trick marks the cell
    # with a "quotePrefix" style under the hood, so reset back to the
    # workbook's plain Normal style right after - callers that want a
    # different format (e.g. the bold header style) re-apply it afterwards
    # via Copy-CellFormat.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

function Copy-CellFormat($srcRange, $dstRange) {
    # Copy formats only (no values) - reuses an existing style index
    # instead of synthesizing a new one.
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 1) Duplicate sheet 2 ("2022-Q3") so its current data survives on its
#    own sheet once sheet 2 itself gets repurposed for 2022-Q4.
# ---------------------------------------------------------------------
$q3Source = $wb.Worksheets.Item(2)
$q3Source.Copy($null, $q3Source)

# ---------------------------------------------------------------------
# 2) Rename sheet 2 to "2022-Q4" and replace its fund table with the new
#    quarter's data.
# ---------------------------------------------------------------------
$q3Source.Name = "2022-Q4"
$q4 = $q3Source

$q4.Cells.Clear()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $q4.Cells.Item(1, $col)
    Set-TextValue $cell $headers[$col - 2]
}

$q4Rows = @(
    , @("012093", "鹏华创新升级混合A", "5.64", "75.32", "3.31", "0.1867", 10)
    , @("014313", "鹏华创新增长一年持有期混合A", "3.47", "60.62", "3.27", "0.1135", 7)
    , @("012094", "鹏华创新升级混合C", "0.17", "75.32", "3.31", "0.0056", 10)
    , @("014314", "鹏华创新增长一年持有期混合C", "0.15", "60.62", "3.27", "0.0049", 7)
)

for ($i = 0; $i -lt $q4Rows.Count; $i++) {
    $r = $i + 2
    $row = $q4Rows[$i]
    $q4.Cells.Item($r, 1).Value = $i

    for ($col = 2; $col -le 7; $col++) {
        $cell = $q4.Cells.Item($r, $col)
        Set-TextValue $cell $row[$col - 2]
    }

    $q4.Cells.Item($r, 8).Value = $row[6]
}

# Re-apply the header / index-column styling (style carried the bold
# centred format, index 2 in the original styles.xml) now that Clear()
# wiped it.
Copy-CellFormat $wb.Worksheets.Item(1).Range("B1") $q4.Range("B1:H1")
Copy-CellFormat $wb.Worksheets.Item(1).Range("A2") $q4.Range("A2:A5")

# ---------------------------------------------------------------------
# 3) Rename the duplicated sheet to "2022-Q3" - its data/format already
#    match the pre-edit sheet, nothing further to change.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(3)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 4) Update the "总计" summary sheet with the new 2022-Q4 row, pushing the
#    existing 2022-Q3 totals down to row 3.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Range("A3").Value = 1
Set-TextValue $summary.Range("B3") "2022-Q3"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.05
Copy-CellFormat $summary.Range("A2") $summary.Range("A3")

$summary.Range("A2").Value = 0
Set-TextValue $summary.Range("B2") "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.31
